# Update the "想去人数" (F) and "最低票价" (G) figures on both the
# "展览" and "全部类型" worksheets, which hold duplicated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F5").Value = 7044
    $ws.Range("F6").Value = 3894
    $ws.Range("F7").Value = 65
    $ws.Range("F8").Value = 158
    $ws.Range("F12").Value = 53
    $ws.Range("F15").Value = 603
    $ws.Range("G15").Value = 58
    $ws.Range("F16").Value = 64
}
